# Content - Add jawfrey coefficient for generating missions
#
# Inserts a new row into the "missionDragonModifiersDefinitions" table
# (Table13303132) for the new dragon sku "dragon_jawfrey", with
# quantityModifier = 4.5 and missionSCRewardMultiplier = 62.
# The new row is inserted right before the existing "dragon_balrog" row
# (which pushes dragon_balrog and dragon_titan down by one row), matching
# the table's descending order by quantityModifier.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tableName = "Table13303132"
$lo = $ws.ListObjects.Item($tableName)

# Row 84 currently holds the "dragon_balrog" entry (5 / 67); insert a new
# blank row above it so it - and everything below it - shifts down by one.
$insertRowNumber = 84
$ws.Rows("$($insertRowNumber):$($insertRowNumber)").Insert()

# Populate the freshly inserted row with the new dragon's data, copying
# the cell styles used by the rest of the table's data rows.
$newRow = $ws.Range("B$($insertRowNumber):E$($insertRowNumber)")
$templateRow = $ws.Range("B$($insertRowNumber + 1):E$($insertRowNumber + 1)")
$newRow.Style = "Normal"
$ws.Range("B$($insertRowNumber)").Value = "<Definition>"
$ws.Range("C$($insertRowNumber)").Value = "dragon_jawfrey"
$ws.Range("D$($insertRowNumber)").Value = 4.5
$ws.Range("E$($insertRowNumber)").Value = 62

$ws.Cells.Item($insertRowNumber, 2).NumberFormat = $ws.Cells.Item($insertRowNumber + 1, 2).NumberFormat
$ws.Cells.Item($insertRowNumber, 3).NumberFormat = $ws.Cells.Item($insertRowNumber + 1, 3).NumberFormat
$ws.Cells.Item($insertRowNumber, 4).NumberFormat = $ws.Cells.Item($insertRowNumber + 1, 4).NumberFormat
$ws.Cells.Item($insertRowNumber, 5).NumberFormat = $ws.Cells.Item($insertRowNumber + 1, 5).NumberFormat

# Copy over formatting (font color, fills, borders, alignment) from the
# row below (the old "dragon_balrog" row, now shifted down) so the new
# row is visually consistent with the rest of the table.
$templateRow.Copy()
$newRow.PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B$($insertRowNumber)").Value = "<Definition>"
$ws.Range("C$($insertRowNumber)").Value = "dragon_jawfrey"
$ws.Range("D$($insertRowNumber)").Value = 4.5
$ws.Range("E$($insertRowNumber)").Value = 62

# The plain row-insert above does not automatically grow the ListObject
# (table) definition, so resize it explicitly to include the new row.
$lo.Resize($ws.Range("B75:E86"))

# The two tables further down the sheet (Difficulty Modifiers, Other
# Modifiers) keep their original size but need to slide down one row
# to follow the rows that got pushed down by the insert above.
$loDifficultyModifiers = $ws.ListObjects.Item("Table1330313234")
$loDifficultyModifiers.Resize($ws.Range("B90:D93"))

$loOtherModifiers = $ws.ListObjects.Item("Table133031323435")
$loOtherModifiers.Resize($ws.Range("B97:D98"))

$excel.CutCopyMode = $false
